# "updated CRUD and added presentation to documentation folder"
#
# Adds a new "phpMyAdmin" / "C" row to the CRUD matrix (between the
# existing "Checkout" row and the thick separator line above the
# footnote), and adds a trailing blank spacer row after the footnote.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$xlPasteFormats = [Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats

# 1) Insert a brand-new row above the current row 10 (the thick
#    separator row). This pushes the separator row 10 -> 11 and the
#    footnote row 11 -> 12, carrying their formatting with them.
$ws.Rows("10:10").Insert()

# 2) Populate the new row 10 with the phpMyAdmin / Create entry, reusing
#    the look of the other table rows (copy formatting from row 8, then
#    set the text so the existing cell styles - incl. the double
#    border/thick rules - are reused instead of new ones being minted).
$ws.Range("A8:E8").Copy()
$ws.Range("A10:E10").PasteSpecial($xlPasteFormats)
$ws.Range("A10").Value = "phpMyAdmin"
$ws.Range("B10").Value = "C"
$ws.Rows(10).RowHeight = 16.5

# 3) The blank cell above the footnote (now row 11) picks up the same
#    style as the footnote text cell.
$ws.Range("A12").Copy()
$ws.Range("A11").PasteSpecial($xlPasteFormats)

# 4) Add a trailing blank row (row 13) below the footnote, styled the
#    same way.
$ws.Range("A12").Copy()
$ws.Range("A13").PasteSpecial($xlPasteFormats)

$excel.CutCopyMode = $false

# Restore the selection to where the user last clicked.
$ws.Range("C19").Select()
